$wb = $excel.ActiveWorkbook

# --- Sheet "Repayment schedule": add column O (value 0) for rows 2-15 ---
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

for ($r = 2; $r -le 15; $r++) {
    # Copy the formatting from the neighboring "N" column cell (same style
    # used throughout the table) onto the new "O" cell before assigning its
    # value, so the new cell reuses the existing shared style instead of
    # Excel minting a brand new cell style.
    $wsSchedule.Range("N$r").Copy()
    $wsSchedule.Range("O$r").PasteSpecial(-4122)
    $wsSchedule.Range("O$r").Value = 0
}

# Update the stored selection for this sheet to the full 16th row.
$wsSchedule.Rows.Item(16).Select()

# --- Sheet "Transactions": renumber the ID column and update selection ---
$wsTransactions = $wb.Worksheets.Item("Transactions")

$wsTransactions.Range("A2").Value = 76
$wsTransactions.Range("A3").Value = 75
$wsTransactions.Range("A4").Value = 74

# Select D2 last so that "Transactions" remains the active sheet/tab, matching
# the original workbook state (tabSelected="1" on this sheet).
$wsTransactions.Range("D2").Select()
